$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-25 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-26 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("73×83=6059", $true, $false, $false, $false, $false, $true, 1, $false, "31×48=1488", 2) | Out-Null
$d.Content.Find.Execute("70×43=3010", $true, $false, $false, $false, $false, $true, 1, $false, "96×79=7584", 2) | Out-Null
$d.Content.Find.Execute("40×52=2080", $true, $false, $false, $false, $false, $true, 1, $false, "46×81=3726", 2) | Out-Null
$d.Content.Find.Execute("56×77=4312", $true, $false, $false, $false, $false, $true, 1, $false, "93×97=9021", 2) | Out-Null
$d.Content.Find.Execute("87×93=8091", $true, $false, $false, $false, $false, $true, 1, $false, "48×28=1344", 2) | Out-Null
$d.Content.Find.Execute("50×30=1500", $true, $false, $false, $false, $false, $true, 1, $false, "45×45=2025", 2) | Out-Null
$d.Content.Find.Execute("76×59=4484", $true, $false, $false, $false, $false, $true, 1, $false, "14×24=336", 2) | Out-Null
$d.Content.Find.Execute("22×60=1320", $true, $false, $false, $false, $false, $true, 1, $false, "90×16=1440", 2) | Out-Null
$d.Content.Find.Execute("49×23=1127", $true, $false, $false, $false, $false, $true, 1, $false, "15×70=1050", 2) | Out-Null
$d.Content.Find.Execute("73×23=1679", $true, $false, $false, $false, $false, $true, 1, $false, "23×61=1403", 2) | Out-Null
$d.Content.Find.Execute("52×55=2860", $true, $false, $false, $false, $false, $true, 1, $false, "27×16=432", 2) | Out-Null
$d.Content.Find.Execute("82×96=7872", $true, $false, $false, $false, $false, $true, 1, $false, "48×45=2160", 2) | Out-Null
$d.Content.Find.Execute("15×53=795", $true, $false, $false, $false, $false, $true, 1, $false, "70×73=5110", 2) | Out-Null
$d.Content.Find.Execute("93×46=4278", $true, $false, $false, $false, $false, $true, 1, $false, "59×22=1298", 2) | Out-Null
$d.Content.Find.Execute("25×26=650", $true, $false, $false, $false, $false, $true, 1, $false, "76×73=5548", 2) | Out-Null
$d.Content.Find.Execute("35×38=1330", $true, $false, $false, $false, $false, $true, 1, $false, "62×72=4464", 2) | Out-Null
$d.Content.Find.Execute("71×11=781", $true, $false, $false, $false, $false, $true, 1, $false, "17×11=187", 2) | Out-Null
$d.Content.Find.Execute("82×55=4510", $true, $false, $false, $false, $false, $true, 1, $false, "91×59=5369", 2) | Out-Null
$d.Content.Find.Execute("91×45=4095", $true, $false, $false, $false, $false, $true, 1, $false, "56×44=2464", 2) | Out-Null
$d.Content.Find.Execute("48×37=1776", $true, $false, $false, $false, $false, $true, 1, $false, "45×62=2790", 2) | Out-Null
$d.Content.Find.Execute("26×44=1144", $true, $false, $false, $false, $false, $true, 1, $false, "66×14=924", 2) | Out-Null
$d.Content.Find.Execute("24×61=1464", $true, $false, $false, $false, $false, $true, 1, $false, "82×38=3116", 2) | Out-Null
$d.Content.Find.Execute("92×88=8096", $true, $false, $false, $false, $false, $true, 1, $false, "93×96=8928", 2) | Out-Null
$d.Content.Find.Execute("57×37=2109", $true, $false, $false, $false, $false, $true, 1, $false, "73×99=7227", 2) | Out-Null
$d.Content.Find.Execute("29×35=1015", $true, $false, $false, $false, $false, $true, 1, $false, "59×61=3599", 2) | Out-Null
